$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46007
$ws.Range("B2").Value = 95.73
$ws.Range("C2").Value = 86
$ws.Range("D2").Value = 81.48
$ws.Range("E2").Value = 76.45999999999999
$ws.Range("F2").Value = 74.5
$ws.Range("G2").Value = 76.13
$ws.Range("H2").Value = 87.63
$ws.Range("I2").Value = 93.17
$ws.Range("J2").Value = 99.23
$ws.Range("K2").Value = 99.06
$ws.Range("L2").Value = 92.69
$ws.Range("M2").Value = 89.45
$ws.Range("N2").Value = 88.5
$ws.Range("O2").Value = 85.90000000000001
$ws.Range("P2").Value = 88.48999999999999
$ws.Range("Q2").Value = 95.70999999999999
$ws.Range("R2").Value = 99.42
$ws.Range("S2").Value = 107.73
$ws.Range("T2").Value = 122.4
$ws.Range("U2").Value = 132.92
$ws.Range("V2").Value = 115.37
$ws.Range("W2").Value = 102.06
$ws.Range("X2").Value = 93.68000000000001
$ws.Range("Y2").Value = 84.95
$ws.Range("Z2").Value = 94.53
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 115.62
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 127.66
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 108.72
$ws.Range("AG2").Value = "1h-23h"
